$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (from AC1, style index "1": bold, bordered,
# centered/top-aligned) onto the three new header cells so they match the
# look of the rest of the header row.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# New header labels
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every data row (2 through 49) gets the same team record: 95 wins, 67
# losses, 0 ties.
for ($row = 2; $row -le 49; $row++) {
    $ws.Cells.Item($row, 30).Value = 95   # column AD
    $ws.Cells.Item($row, 31).Value = 67   # column AE
    $ws.Cells.Item($row, 32).Value = 0    # column AF
}

Write-Output "Added Wins/Losses/Ties columns (AD:AF) to rows 1-49"
